$d = $word.ActiveDocument

# --- Paragraph 1: title -> Heading1 styled "Rent Ledger Summary" ---
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Rent Ledger Summary</w:t></w:r></w:p>
'@
[void]$d.Paragraphs(1).Range.InsertXML($p1xml)

# --- Paragraph 2: "Date: 2026-01-21" -> "Property: {{PROPERTY_NAME}}" ---
$d.Paragraphs(2).Range.Text = "Property: {{PROPERTY_NAME}}"

# --- Paragraph 3: "Not legal advice..." + <w:br/> -> "Period: {{PERIOD_START}} to {{PERIOD_END}}" ---
$d.Paragraphs(3).Range.Text = "Period: {{PERIOD_START}} to {{PERIOD_END}}"

# --- Paragraph 4: "Date | Type | ..." -> 6x2 table ---
$tblXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblW w:type="auto" w:w="0"/>
    <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1440"/>
    <w:gridCol w:w="1440"/>
    <w:gridCol w:w="1440"/>
    <w:gridCol w:w="1440"/>
    <w:gridCol w:w="1440"/>
    <w:gridCol w:w="1440"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>Date</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>Tenant</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>Unit</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>Charge Type</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>Amount</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>Balance</w:t></w:r></w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>{{DATE}}</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>{{TENANT}}</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>{{UNIT}}</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>{{CHARGE_TYPE}}</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>{{AMOUNT}}</w:t></w:r></w:p>
    </w:tc>
    <w:tc>
      <w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr>
      <w:p><w:r><w:t>{{BALANCE}}</w:t></w:r></w:p>
    </w:tc>
  </w:tr>
</w:tbl>
'@
[void]$d.Paragraphs(4).Range.InsertXML($tblXml)

Write-Host "Edit complete"
